$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1592.5714
$ws.Range("I4").Value = 189.2
$ws.Range("K4").Value = 189.2
$ws.Range("M4").Value = -75.19999999999999
$ws.Range("H5").Value = 192.16667
$ws.Range("I5").Value = 192.16667
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 192.16667
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -77.16667000000001
$ws.Range("N5").ClearContents()
$ws.Range("H70").Value = 4601
$ws.Range("J70").Value = 4901.0835
$ws.Range("L70").Value = 14703.2505
$ws.Range("N70").Value = -15243.2505
$ws.Range("H73").Value = 4601
$ws.Range("J73").Value = 4901.0835
$ws.Range("L73").Value = 14703.2505
$ws.Range("N73").Value = -16575.2505
$ws.Range("H98").Value = 1044.875
$ws.Range("J98").Value = 1051.25
$ws.Range("L98").Value = 1051.25
$ws.Range("N98").Value = -4047.25
$ws.Range("H101").Value = 342.16666
$ws.Range("I101").Value = 199
$ws.Range("J101").Value = 370.8
$ws.Range("K101").Value = 597
$ws.Range("L101").Value = 1112.4
$ws.Range("M101").Value = 1025
$ws.Range("N101").Value = -4356.4
$ws.Range("H103").Value = 3010.1428
$ws.Range("J103").Value = 2047.5
$ws.Range("L103").Value = 6142.5
$ws.Range("N103").Value = -7314.5
$ws.Range("H113").Value = 4999
$ws.Range("J113").Value = 4998
$ws.Range("L113").Value = 4998
$ws.Range("N113").Value = -11506
$ws.Range("H122").Value = 1044.875
$ws.Range("J122").Value = 1051.25
$ws.Range("L122").Value = 3153.75
$ws.Range("N122").Value = -8053.75
$ws.Range("H138").Value = 3219.1667
$ws.Range("I138").Value = 1968.2858
$ws.Range("J138").Value = 4015.182
$ws.Range("K138").Value = 5904.857400000001
$ws.Range("L138").Value = 12045.546
$ws.Range("M138").Value = -764.8574000000008
$ws.Range("N138").Value = -22325.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11786.5
$ws.Range("I2").Value = 8560
$ws.Range("J2").Value = 15013
$ws.Range("K2").Value = 8560
$ws.Range("L2").Value = 15013
$ws.Range("M2").Value = -8447
$ws.Range("N2").Value = -15239
$ws.Range("H32").Value = 9098567
$ws.Range("J32").Value = 20013600
$ws.Range("L32").Value = 20013600
$ws.Range("N32").Value = -20014174
$ws.Range("H61").Value = 1330.4375
$ws.Range("I61").Value = 1359.5
$ws.Range("J61").Value = 1127
$ws.Range("K61").Value = 1359.5
$ws.Range("L61").Value = 1127
$ws.Range("M61").Value = -1147.5
$ws.Range("N61").Value = -1551
$ws.Range("H74").Value = 5861.9375
$ws.Range("I74").Value = 6141.4165
$ws.Range("J74").Value = 5023.5
$ws.Range("K74").Value = 6141.4165
$ws.Range("L74").Value = 5023.5
$ws.Range("M74").Value = -5267.4165
$ws.Range("N74").Value = -6771.5
$ws.Range("H77").Value = 5861.9375
$ws.Range("I77").Value = 6141.4165
$ws.Range("J77").Value = 5023.5
$ws.Range("K77").Value = 30707.0825
$ws.Range("L77").Value = 25117.5
$ws.Range("M77").Value = -26339.0825
$ws.Range("N77").Value = -33853.5
$ws.Range("H110").Value = 2291.2666
$ws.Range("J110").Value = 2903.5715
$ws.Range("L110").Value = 2903.5715
$ws.Range("N110").Value = -6993.5715
$ws.Range("H111").Value = 34999
$ws.Range("J111").Value = 34999
$ws.Range("L111").Value = 34999
$ws.Range("N111").Value = -43179
$ws.Range("H116").Value = 11786.5
$ws.Range("I116").Value = 8560
$ws.Range("J116").Value = 15013
$ws.Range("K116").Value = 8560
$ws.Range("L116").Value = 15013
$ws.Range("M116").Value = -6266
$ws.Range("N116").Value = -19601
$ws.Range("H122").Value = 2550.889
$ws.Range("I122").Value = 2542
$ws.Range("K122").Value = 7626
$ws.Range("M122").Value = -5176
$ws.Range("H132").Value = 2384.889
$ws.Range("I132").Value = 2486.652
$ws.Range("K132").Value = 7459.956
$ws.Range("M132").Value = -4929.956
$ws.Range("H136").Value = 1330.4375
$ws.Range("I136").Value = 1359.5
$ws.Range("J136").Value = 1127
$ws.Range("K136").Value = 4078.5
$ws.Range("L136").Value = 3381
$ws.Range("M136").Value = -1528.5
$ws.Range("N136").Value = -8481

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11786.5
$ws.Range("I3").Value = 8560
$ws.Range("J3").Value = 15013
$ws.Range("K3").Value = 8560
$ws.Range("L3").Value = 15013
$ws.Range("M3").Value = -8446
$ws.Range("N3").Value = -15241
$ws.Range("H134").Value = 2284.75
$ws.Range("I134").Value = 1974.3636
$ws.Range("J134").Value = 5699
$ws.Range("K134").Value = 5923.0908
$ws.Range("L134").Value = 17097
$ws.Range("M134").Value = -3388.0908
$ws.Range("N134").Value = -22167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -387
$ws.Range("N2").ClearContents()
$ws.Range("H7").Value = 85.8
$ws.Range("J7").Value = 164.66667
$ws.Range("L7").Value = 164.66667
$ws.Range("N7").Value = -390.66667
$ws.Range("H31").Value = 5693.0293
$ws.Range("I31").Value = 3378.7
$ws.Range("K31").Value = 3378.7
$ws.Range("M31").Value = -3083.7
$ws.Range("H34").Value = 5693.0293
$ws.Range("I34").Value = 3378.7
$ws.Range("K34").Value = 3378.7
$ws.Range("M34").Value = -3176.7
$ws.Range("H86").Value = 4999.5
$ws.Range("I86").Value = 4999.5
$ws.Range("K86").Value = 4999.5
$ws.Range("M86").Value = -3876.5
$ws.Range("H89").Value = 4999.5
$ws.Range("I89").Value = 4999.5
$ws.Range("K89").Value = 24997.5
$ws.Range("M89").Value = -19381.5
$ws.Range("H94").Value = 4045.9
$ws.Range("I94").Value = 1226.5
$ws.Range("K94").Value = 1226.5
$ws.Range("M94").Value = -775.5
$ws.Range("H107").Value = 1127.625
$ws.Range("I107").Value = 915.36365
$ws.Range("K107").Value = 915.36365
$ws.Range("M107").Value = 1004.63635
$ws.Range("H122").Value = 802.9
$ws.Range("I122").Value = 689.375
$ws.Range("K122").Value = 2068.125
$ws.Range("M122").Value = 381.875
$ws.Range("H132").Value = 4278.3335
$ws.Range("I132").Value = 3917.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 11752.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -9222.5
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 2224.625
$ws.Range("I134").Value = 2224.625
$ws.Range("K134").Value = 6673.875
$ws.Range("M134").Value = -4138.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4952.2085
$ws.Range("I80").Value = 4999.9443
$ws.Range("J80").Value = 4809
$ws.Range("K80").Value = 14999.8329
$ws.Range("L80").Value = 14427
$ws.Range("M80").Value = -14063.8329
$ws.Range("N80").Value = -16299
$ws.Range("H83").Value = 4952.2085
$ws.Range("I83").Value = 4999.9443
$ws.Range("J83").Value = 4809
$ws.Range("K83").Value = 44999.4987
$ws.Range("L83").Value = 43281
$ws.Range("M83").Value = -40319.4987
$ws.Range("N83").Value = -52641
$ws.Range("H103").Value = 2268.6365
$ws.Range("J103").Value = 2268.6365
$ws.Range("L103").Value = 6805.9095
$ws.Range("N103").Value = -8563.9095
$ws.Range("H114").Value = 3499.6667
$ws.Range("I114").Value = 4999.5
$ws.Range("K114").Value = 14998.5
$ws.Range("M114").Value = -11744.5
$ws.Range("H122").Value = 823.2
$ws.Range("I122").Value = 840.7143
$ws.Range("J122").Value = 782.3333
$ws.Range("K122").Value = 7566.428699999999
$ws.Range("L122").Value = 7040.9997
$ws.Range("M122").Value = -5116.428699999999
$ws.Range("N122").Value = -11940.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1467.7858
$ws.Range("I102").Value = 1467.7858
$ws.Range("K102").Value = 1467.7858
$ws.Range("M102").Value = 154.2141999999999
$ws.Range("H132").Value = 76048.64
$ws.Range("I132").Value = 128254
$ws.Range("K132").Value = 384762
$ws.Range("M132").Value = -382232

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1271.1364
$ws.Range("I93").Value = 1203.2778
$ws.Range("J93").Value = 1576.5
$ws.Range("K93").Value = 1203.2778
$ws.Range("L93").Value = 1576.5
$ws.Range("M93").Value = 44.72219999999993
$ws.Range("N93").Value = -4072.5
$ws.Range("H100").Value = 10000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8228
$ws.Range("I62").Value = 5513.25
$ws.Range("J62").Value = 10399.8
$ws.Range("K62").Value = 5513.25
$ws.Range("L62").Value = 10399.8
$ws.Range("M62").Value = -4889.25
$ws.Range("N62").Value = -11647.8
$ws.Range("H65").Value = 8228
$ws.Range("I65").Value = 5513.25
$ws.Range("J65").Value = 10399.8
$ws.Range("K65").Value = 27566.25
$ws.Range("L65").Value = 51999
$ws.Range("M65").Value = -24446.25
$ws.Range("N65").Value = -58239
$ws.Range("H93").Value = 33389
$ws.Range("J93").Value = 33389
$ws.Range("L93").Value = 33389
$ws.Range("N93").Value = -38381
$ws.Range("H96").Value = 1081.8334
$ws.Range("I96").Value = 1262
$ws.Range("J96").Value = 829.6
$ws.Range("K96").Value = 1262
$ws.Range("L96").Value = 829.6
$ws.Range("M96").Value = 111
$ws.Range("N96").Value = -3575.6
